$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.920.13'
$ws.Range("E2").Value = '  +3.88%  '
$ws.Range("D3").Value = '2.655.14'
$ws.Range("E3").Value = '  +6.30%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.60'
$ws.Range("E5").Value = '  +7.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '328.06'
$ws.Range("E6").Value = '  +2.98%  '
$ws.Range("E7").Value = '  +1.83%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.554'
$ws.Range("E9").Value = '  +3.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.25'
$ws.Range("E10").Value = '  +6.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.18'
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0823'
$ws.Range("E12").Value = '  +2.79%  '
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.36'
$ws.Range("E14").Value = '  +4.61%  '
$ws.Range("D15").Value = '3.067.96'
$ws.Range("E15").Value = '  +6.13%  '
$ws.Range("D16").Value = '2.668.83'
$ws.Range("E16").Value = '  +6.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.874'
$ws.Range("E17").Value = '  +5.46%  '
$ws.Range("D18").Value = '49.851.94'
$ws.Range("E18").Value = '  +4.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.17'
$ws.Range("E19").Value = '  +2.37%  '
$ws.Range("E20").Value = '  +2.24%  '
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("D22").Value = '0.0₃0959'
$ws.Range("E22").Value = '  +3.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.29'
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '276.82'
$ws.Range("E24").Value = '  +1.83%  '
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("E26").Value = '  +4.75%  '
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  +3.01%  '
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.98'
$ws.Range("E30").Value = '  +4.07%  '
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.25'
$ws.Range("E32").Value = '  +1.83%  '
$ws.Range("E33").Value = '  +2.93%  '
$ws.Range("E34").Value = '  +3.13%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0803'
$ws.Range("E35").Value = '  +3.67%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  +7.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.78'
$ws.Range("E38").Value = '  +4.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.07'
$ws.Range("E39").Value = '  +7.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.63'
$ws.Range("E40").Value = '  +4.67%  '
$ws.Range("E41").Value = '  +2.16%  '
$ws.Range("E42").Value = '  +1.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.28'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0316'
$ws.Range("E44").Value = '  +4.57%  '
$ws.Range("E45").Value = '  +5.63%  '
$ws.Range("D46").Value = '2.070.15'
$ws.Range("E46").Value = '  +3.69%  '
$ws.Range("E47").Value = '  +12.71%  '
$ws.Range("E48").Value = '  +5.73%  '
$ws.Range("E49").Value = '  +3.39%  '
$ws.Range("E50").Value = '  +4.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '82.08'
$ws.Range("E51").Value = '  +4.41%  '
